$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 210.3866666666667
$ws.Range("C2").Value = 201.8933333333333
$ws.Range("D2").Value = 282.92
$ws.Range("E2").Value = 240.1733333333333
$ws.Range("F2").Value = 252.5333333333333
$ws.Range("G2").Value = 257.3199999999999
$ws.Range("H2").Value = 488.4533333333334
$ws.Range("I2").Value = 466.96
$ws.Range("J2").Value = 482.64
$ws.Range("K2").Value = 423.7066666666666
$ws.Range("L2").Value = 385.4133333333333
$ws.Range("M2").Value = 382.3333333333333
$ws.Range("B3").Value = 150.1066666666667
$ws.Range("C3").Value = 210.3066666666666
$ws.Range("D3").Value = 259.7333333333333
$ws.Range("E3").Value = 261.9066666666666
$ws.Range("F3").Value = 257.4533333333333
$ws.Range("G3").Value = 264.04
$ws.Range("H3").Value = 493.2666666666667
$ws.Range("I3").Value = 468.6
$ws.Range("J3").Value = 510.3200000000001
$ws.Range("K3").Value = 439.5333333333333
$ws.Range("L3").Value = 435.1466666666666
$ws.Range("M3").Value = 368.5466666666667
$ws.Range("B4").Value = 99.13333333333333
$ws.Range("C4").Value = 118.4
$ws.Range("D4").Value = 128.0533333333333
$ws.Range("E4").Value = 134.12
$ws.Range("F4").Value = 137.9333333333333
$ws.Range("G4").Value = 144.2533333333333
$ws.Range("H4").Value = 318.5066666666667
$ws.Range("I4").Value = 304.16
$ws.Range("J4").Value = 289.24
$ws.Range("K4").Value = 279.04
$ws.Range("L4").Value = 270.2666666666667
$ws.Range("M4").Value = 229.68
$ws.Range("B5").Value = 47.66666666666666
$ws.Range("C5").Value = 53.65333333333334
$ws.Range("D5").Value = 76.21333333333332
$ws.Range("E5").Value = 57.12
$ws.Range("F5").Value = 81.78666666666666
$ws.Range("G5").Value = 69.53333333333333
$ws.Range("H5").Value = 257.1733333333333
$ws.Range("I5").Value = 232.28
$ws.Range("J5").Value = 255.48
$ws.Range("K5").Value = 214.0133333333333
$ws.Range("L5").Value = 194.32
$ws.Range("M5").Value = 168.0133333333333
$ws.Range("B6").Value = 42.25333333333333
$ws.Range("C6").Value = 44.78666666666666
$ws.Range("D6").Value = 49.22666666666667
$ws.Range("E6").Value = 48.53333333333332
$ws.Range("F6").Value = 49.92
$ws.Range("G6").Value = 39.33333333333333
$ws.Range("H6").Value = 194.76
$ws.Range("I6").Value = 183.36
$ws.Range("J6").Value = 175.8
$ws.Range("K6").Value = 153.2933333333333
$ws.Range("L6").Value = 139.8133333333333
$ws.Range("M6").Value = 108.5866666666667
$ws.Range("H7").Value = 62.65333333333334
$ws.Range("I7").Value = 58.38666666666666
$ws.Range("J7").Value = 58.56
$ws.Range("K7").Value = 37.12
$ws.Range("B8").Value = 242.2266666666667
$ws.Range("C8").Value = 271.3066666666666
$ws.Range("D8").Value = 254.0133333333333
$ws.Range("E8").Value = 225.68
$ws.Range("F8").Value = 219.4533333333333
$ws.Range("G8").Value = 234.96
$ws.Range("H8").Value = 400.0933333333334
$ws.Range("I8").Value = 361.5466666666667
$ws.Range("J8").Value = 364.0933333333333
$ws.Range("K8").Value = 340.6666666666666
$ws.Range("L8").Value = 331.8933333333333
$ws.Range("M8").Value = 342.88
$ws.Range("B9").Value = 127.24
$ws.Range("C9").Value = 240.88
$ws.Range("D9").Value = 218.64
$ws.Range("E9").Value = 192.76
$ws.Range("F9").Value = 187.3066666666667
$ws.Range("G9").Value = 183.68
$ws.Range("H9").Value = 364.1866666666667
$ws.Range("I9").Value = 357.48
$ws.Range("J9").Value = 352.1333333333333
$ws.Range("K9").Value = 330.5599999999999
$ws.Range("L9").Value = 341.5466666666666
$ws.Range("M9").Value = 344.4533333333333
$ws.Range("B10").Value = 25.25333333333334
$ws.Range("C10").Value = 170.3466666666666
$ws.Range("D10").Value = 141.4933333333333
$ws.Range("E10").Value = 119.7466666666667
$ws.Range("F10").Value = 129.0266666666666
$ws.Range("G10").Value = 117.2133333333333
$ws.Range("H10").Value = 311.5333333333333
$ws.Range("I10").Value = 276.3466666666667
$ws.Range("J10").Value = 285.6799999999999
$ws.Range("K10").Value = 268.04
$ws.Range("L10").Value = 251.76
$ws.Range("M10").Value = 249.5333333333333
$ws.Range("C11").Value = 45.48
$ws.Range("D11").Value = 93.37333333333333
$ws.Range("E11").Value = 83.90666666666667
$ws.Range("F11").Value = 59.64
$ws.Range("G11").Value = 32.53333333333333
$ws.Range("H11").Value = 269
$ws.Range("I11").Value = 247.24
$ws.Range("J11").Value = 239.2666666666667
$ws.Range("K11").Value = 239.8133333333333
$ws.Range("L11").Value = 237.52
$ws.Range("M11").Value = 205.1066666666667
$ws.Range("H12").Value = 244.16
$ws.Range("I12").Value = 212.9466666666667
$ws.Range("J12").Value = 188.6266666666667
$ws.Range("K12").Value = 183.6533333333333
$ws.Range("L12").Value = 187.4133333333333
$ws.Range("M12").Value = 105.7866666666667
$ws.Range("H13").Value = 178.2133333333333
$ws.Range("I13").Value = 143.9866666666667
$ws.Range("J13").Value = 99.75999999999999
$ws.Range("K13").Value = 72.82666666666667
$ws.Range("L13").Value = 40.88
$ws.Range("C14").Value = 154.72
$ws.Range("D14").Value = 114.92
$ws.Range("E14").Value = 124.56
$ws.Range("F14").Value = 120.04
$ws.Range("G14").Value = 112.8933333333333
$ws.Range("H14").Value = 241.84
$ws.Range("I14").Value = 209.72
$ws.Range("J14").Value = 198.4933333333333
$ws.Range("K14").Value = 193.4533333333333
$ws.Range("L14").Value = 216.7733333333333
$ws.Range("M14").Value = 214.92
$ws.Range("C15").Value = 77.73333333333332
$ws.Range("D15").Value = 88.46666666666667
$ws.Range("E15").Value = 75.31999999999999
$ws.Range("F15").Value = 84.19999999999999
$ws.Range("G15").Value = 96.88
$ws.Range("H15").Value = 253.9333333333333
$ws.Range("I15").Value = 226.2266666666667
$ws.Range("J15").Value = 211.9333333333333
$ws.Range("K15").Value = 211.7333333333333
$ws.Range("L15").Value = 182.5466666666667
$ws.Range("M15").Value = 239.12
$ws.Range("C16").Value = 16.53333333333333
$ws.Range("D16").Value = 21.33333333333333
$ws.Range("E16").Value = 20.97333333333334
$ws.Range("F16").Value = 20.48
$ws.Range("G16").Value = 16.49333333333333
$ws.Range("H16").Value = 159.4666666666666
$ws.Range("I16").Value = 143.6133333333333
$ws.Range("J16").Value = 153.4
$ws.Range("K16").Value = 120.6933333333333
$ws.Range("L16").Value = 124.5066666666667
$ws.Range("M16").Value = 121.2133333333333
$ws.Range("D17").Value = 0
$ws.Range("H17").Value = 161.8933333333333
$ws.Range("I17").Value = 153.2533333333333
$ws.Range("J17").Value = 151.9466666666667
$ws.Range("K17").Value = 125.5733333333333
$ws.Range("L17").Value = 120.5733333333333
$ws.Range("M17").Value = 118.9733333333333
$ws.Range("H18").Value = 122.7466666666667
$ws.Range("I18").Value = 129.76
$ws.Range("J18").Value = 155.6933333333333
$ws.Range("K18").Value = 109.4
$ws.Range("L18").Value = 98.31999999999999
$ws.Range("M18").Value = 90.8
$ws.Range("H19").Value = 41.65333333333334
$ws.Range("I19").Value = 35.25333333333334
$ws.Range("J19").Value = 47.01333333333334
$ws.Range("K19").Value = 40.54666666666667
$ws.Range("L19").Value = 40.94666666666666
$ws.Range("M19").Value = 60.13333333333333
$ws.Range("C20").Value = 107.1466666666666
$ws.Range("D20").Value = 220.3066666666666
$ws.Range("E20").Value = 159.5333333333333
$ws.Range("F20").Value = 199.0133333333333
$ws.Range("G20").Value = 191.52
$ws.Range("H20").Value = 433.7333333333333
$ws.Range("I20").Value = 435.5333333333333
$ws.Range("J20").Value = 497.52
$ws.Range("K20").Value = 501.9466666666666
$ws.Range("L20").Value = 499.8933333333333
$ws.Range("M20").Value = 489.3466666666667
$ws.Range("C21").Value = 54.73333333333333
$ws.Range("D21").Value = 231.9066666666666
$ws.Range("E21").Value = 246.48
$ws.Range("F21").Value = 286.96
$ws.Range("G21").Value = 355.0666666666666
$ws.Range("H21").Value = 539.0933333333334
$ws.Range("I21").Value = 519.4
$ws.Range("J21").Value = 514.1866666666666
$ws.Range("K21").Value = 521.3466666666666
$ws.Range("L21").Value = 578.6133333333332
$ws.Range("M21").Value = 564.5866666666667
$ws.Range("B22").Value = 83.94666666666666
$ws.Range("C22").Value = 39.2
$ws.Range("D22").Value = 174.68
$ws.Range("E22").Value = 212.7466666666667
$ws.Range("F22").Value = 256.6666666666666
$ws.Range("G22").Value = 244.3066666666667
$ws.Range("H22").Value = 447.6933333333333
$ws.Range("I22").Value = 315.1333333333333
$ws.Range("J22").Value = 349.9733333333334
$ws.Range("K22").Value = 282.7733333333333
$ws.Range("L22").Value = 259.36
$ws.Range("M22").Value = 437.3466666666666
$ws.Range("C23").Value = 21.38666666666666
$ws.Range("D23").Value = 11.2
$ws.Range("E23").Value = 42.79999999999999
$ws.Range("F23").Value = 54.86666666666666
$ws.Range("G23").Value = 154.7466666666667
$ws.Range("H23").Value = 448.96
$ws.Range("I23").Value = 274.9733333333333
$ws.Range("J23").Value = 506.5866666666666
$ws.Range("K23").Value = 289.1066666666667
$ws.Range("L23").Value = 333.76
$ws.Range("M23").Value = 601.7733333333333
